$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 14, shifting the SkillType block down
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the ATTACK attribute value
$ws.Cells.Item(14, 7).Value = "ATTACK"
$ws.Cells.Item(14, 9).Value = 4

# Move the selection to I16, matching the new last populated cell
$null = $ws.Range("I16").Select()
